$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 19: "/*/d1:notificationEF/oos:id" -> "/*/d1:notificationXX/oos:id"
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = "/*/d1:notificationXX/oos:id"

# ---------------------------------------------------------------------------
# 2) Row 20: "/*/d1:notificationEF/oos:notificationNumber" -> "/*/d1:notificationXX/oos:notificationNumber"
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = "/*/d1:notificationXX/oos:notificationNumber"

# ---------------------------------------------------------------------------
# 3) Clear the old "contracts_finances_<region>" / "contracts_<region>" rows
#    that used to live at 22:23 - they get relocated further down (31:32).
# ---------------------------------------------------------------------------
$ws.Range("A22:G23").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 4) New notification rows 21, 23-26 (values) plus matching style for E:G.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "notifications"
$ws.Range("B21").Value = "notifications_<region>"
$ws.Range("C21").Value = "/*/d1:notificationXX/oos:versionNumber"
$ws.Range("D21").Value = "oos:versionNumber"
$ws.Range("E21").Value = '"1"'
$ws.Range("F21").Value = "integer"
$ws.Range("G21").Value = "NotificationVersionNumber"

$ws.Range("A22").Value = "notifications"
$ws.Range("B22").Value = "notifications_<region>"
$ws.Range("C22").Value = "/*/d1:notificationXX/oos:publishDate"
$ws.Range("D22").Value = "oos:publishDate"
$ws.Range("E22").Value = '"2011-01-21T15:31:06"'
$ws.Range("F22").Value = "date"
$ws.Range("G22").Value = "NotificationPublishDate"

$ws.Range("A23").Value = "notifications"
$ws.Range("B23").Value = "notifications_<region>"
$ws.Range("C23").Value = "/*/d1:notificationXX/oos:placingWay/oos:name"
$ws.Range("D23").Value = "oos:placingWay/oos:name"
$ws.Range("E23").Value = '"Открытый конкурс"'
$ws.Range("F23").Value = "factor"
$ws.Range("G23").Value = "NotificationPlacingWayName"

$ws.Range("A24").Value = "notifications"
$ws.Range("B24").Value = "notifications_<region>"
$ws.Range("C24").Value = "/*/d1:notificationXX/oos:orderName"
$ws.Range("D24").Value = "oos:orderName"
$ws.Range("E24").Value = '"Обязательное страхование авто гражданской ответственности владельцев транспортных средств"'
$ws.Range("F24").Value = "character"
$ws.Range("G24").Value = "NotificationOrderName"

$ws.Range("A25").Value = "notifications"
$ws.Range("B25").Value = "notifications_<region>"
$ws.Range("C25").Value = "/*/d1:notificationXX/oos:order/oos:placer/oos:regNum"
$ws.Range("D25").Value = "oos:order/oos:placer/oos:regNum"
$ws.Range("E25").Value = '"01761000012"'
$ws.Range("F25").Value = "factor"
$ws.Range("G25").Value = "NotificationOrderPlacerRegNum"

$ws.Range("A26").Value = "notifications"
$ws.Range("B26").Value = "notifications_<region>"
$ws.Range("C26").Value = "/*/d1:notificationXX/oos:order/oos:placer/oos:fullName"
$ws.Range("D26").Value = "oos:order/oos:placer/oos:fullName"
$ws.Range("E26").Value = '"Управление Федеральной налоговой службы по Республике Адыгея"'
$ws.Range("F26").Value = "factor"
$ws.Range("G26").Value = "NotificationOrderPlacerFullName"

# Give the new rows (21-26) the same "applied number format" style that
# columns E:G otherwise carry throughout this sheet (col E already inherits
# it from the column default, but F/G need it applied explicitly).
$ws.Range("E21:G26").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 5) Three blank spacer rows (27-29), only F/G formatted, no value.
# ---------------------------------------------------------------------------
$ws.Range("F27").NumberFormat = "General"
$ws.Range("G27").NumberFormat = "General"
$ws.Range("F28").NumberFormat = "General"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("F29").NumberFormat = "General"
$ws.Range("G29").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 6) Row 31 / 32: the relocated "contracts_finances_<region>" /
#    "contracts_<region>" rows (previously at 22/23).
# ---------------------------------------------------------------------------
$ws.Range("B31").Value = "contracts_finances_<region>"
$ws.Range("C31").Value = "/d1:contract/oos:finances/oos:financeSource"
$ws.Range("E31").Value = '"Российская Федерация"'
$ws.Range("F31").Value = "factor"
$ws.Range("G31").Value = "ContractFinance"

$ws.Range("B32").Value = "contracts_<region>"
$ws.Range("C32").Value = "/*/d1:contract/oos:foundation/oos:singleCustomer"
$ws.Range("E32").Value = '"true"'
$ws.Range("F32").Value = "factor"
$ws.Range("G32").Value = "ContractSingleCustomer"

$ws.Range("E31:G32").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 7) Row 34: new "notifications_lots_<region>" summary row.
# ---------------------------------------------------------------------------
$ws.Range("B34").Value = "notifications_lots_<region>"
$ws.Range("C34").Value = "(Many lots to one notification)"

# ---------------------------------------------------------------------------
# 8) Update the selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A26").Select() | Out-Null
